# Auto-generated edit script applying numeric corrections to the
# Pandaemonium_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2312.818
$ws.Range("I113").Value = 1626.25
$ws.Range("J113").Value = 2705.1428
$ws.Range("K113").Value = 1626.25
$ws.Range("L113").Value = 2705.1428
$ws.Range("M113").Value = 1627.75
$ws.Range("N113").Value = -9213.1428

$ws.Range("H116").Value = 2124.3103
$ws.Range("I116").Value = 1940.6666
$ws.Range("K116").Value = 1940.6666
$ws.Range("M116").Value = 1501.3334

$ws.Range("H132").Value = 1991.9166
$ws.Range("I132").Value = 1947.826
$ws.Range("J132").Value = 3006
$ws.Range("K132").Value = 5843.478
$ws.Range("L132").Value = 9018
$ws.Range("M132").Value = -3313.478
$ws.Range("N132").Value = -14078

$ws.Range("H138").Value = 21049.588
$ws.Range("I138").Value = 2819.2
$ws.Range("J138").Value = 47093
$ws.Range("K138").Value = 8457.599999999999
$ws.Range("L138").Value = 141279
$ws.Range("M138").Value = -3317.599999999999
$ws.Range("N138").Value = -151559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1961.5385
$ws.Range("I45").Value = 1961.5385
$ws.Range("K45").Value = 1961.5385
$ws.Range("M45").Value = -1584.5385

$ws.Range("H61").Value = 9263.909
$ws.Range("I61").Value = 8657.5
$ws.Range("J61").Value = 11516.286
$ws.Range("K61").Value = 8657.5
$ws.Range("L61").Value = 11516.286
$ws.Range("M61").Value = -8445.5
$ws.Range("N61").Value = -11940.286

$ws.Range("H74").Value = 2333.8
$ws.Range("I74").Value = 1937.0454
$ws.Range("J74").Value = 3424.875
$ws.Range("K74").Value = 1937.0454
$ws.Range("L74").Value = 3424.875
$ws.Range("M74").Value = -1063.0454
$ws.Range("N74").Value = -5172.875

$ws.Range("H77").Value = 2333.8
$ws.Range("I77").Value = 1937.0454
$ws.Range("J77").Value = 3424.875
$ws.Range("K77").Value = 9685.226999999999
$ws.Range("L77").Value = 17124.375
$ws.Range("M77").Value = -5317.226999999999
$ws.Range("N77").Value = -25860.375

$ws.Range("H136").Value = 9263.909
$ws.Range("I136").Value = 8657.5
$ws.Range("J136").Value = 11516.286
$ws.Range("K136").Value = 25972.5
$ws.Range("L136").Value = 34548.858
$ws.Range("M136").Value = -23422.5
$ws.Range("N136").Value = -39648.858

$ws.Range("H140").Value = 46877.5
$ws.Range("J140").Value = 46877.5
$ws.Range("L140").Value = 46877.5
$ws.Range("N140").Value = -57237.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6210.8237
$ws.Range("I105").Value = 6579.6
$ws.Range("J105").Value = 5684
$ws.Range("K105").Value = 6579.6
$ws.Range("L105").Value = 5684
$ws.Range("M105").Value = -4832.6
$ws.Range("N105").Value = -9178

$ws.Range("H134").Value = 3711.8
$ws.Range("I134").Value = 3496.3333
$ws.Range("J134").Value = 3855.4443
$ws.Range("K134").Value = 10488.9999
$ws.Range("L134").Value = 11566.3329
$ws.Range("M134").Value = -7953.999899999999
$ws.Range("N134").Value = -16636.3329

$ws.Range("H138").Value = 49797.5
$ws.Range("J138").Value = 49797.5
$ws.Range("L138").Value = 49797.5
$ws.Range("N138").Value = -60077.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8423.691999999999
$ws.Range("I31").Value = 7407.25
$ws.Range("J31").Value = 10050
$ws.Range("K31").Value = 7407.25
$ws.Range("L31").Value = 10050
$ws.Range("M31").Value = -7112.25
$ws.Range("N31").Value = -10640

$ws.Range("H34").Value = 8423.691999999999
$ws.Range("I34").Value = 7407.25
$ws.Range("J34").Value = 10050
$ws.Range("K34").Value = 7407.25
$ws.Range("L34").Value = 10050
$ws.Range("M34").Value = -7205.25
$ws.Range("N34").Value = -10454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 11912260
$ws.Range("I5").Value = 351.63635
$ws.Range("K5").Value = 1054.90905
$ws.Range("M5").Value = -942.90905

$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H113").Value = 710.01697
$ws.Range("I113").Value = 713.95746
$ws.Range("J113").Value = 694.5833
$ws.Range("K113").Value = 2141.87238
$ws.Range("L113").Value = 2083.7499
$ws.Range("M113").Value = 28.12762000000021
$ws.Range("N113").Value = -6423.7499

$ws.Range("H122").Value = 1205.5385
$ws.Range("I122").Value = 242
$ws.Range("K122").Value = 2178
$ws.Range("M122").Value = 272

$ws.Range("H131").Value = 20699.438
$ws.Range("J131").Value = 29636.697
$ws.Range("L131").Value = 88910.091
$ws.Range("N131").Value = -98990.091

$ws.Range("H135").Value = 11912260
$ws.Range("I135").Value = 351.63635
$ws.Range("K135").Value = 3164.72715
$ws.Range("M135").Value = -629.7271499999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H70").Value = 5536.2163
$ws.Range("I70").Value = 5345.6
$ws.Range("J70").Value = 5933.3335
$ws.Range("K70").Value = 5345.6
$ws.Range("L70").Value = 5933.3335
$ws.Range("M70").Value = -5075.6
$ws.Range("N70").Value = -6473.3335

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H73").Value = 5536.2163
$ws.Range("I73").Value = 5345.6
$ws.Range("J73").Value = 5933.3335
$ws.Range("K73").Value = 5345.6
$ws.Range("L73").Value = 5933.3335
$ws.Range("M73").Value = -4409.6
$ws.Range("N73").Value = -7805.3335

$ws.Range("H126").Value = 1713.7307
$ws.Range("I126").Value = 1683.8334
$ws.Range("J126").Value = 1722.7
$ws.Range("K126").Value = 5051.5002
$ws.Range("L126").Value = 5168.1
$ws.Range("M126").Value = -2581.5002
$ws.Range("N126").Value = -10108.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 493.33334
$ws.Range("I22").Value = 499
$ws.Range("J22").Value = 490.5
$ws.Range("K22").Value = 499
$ws.Range("L22").Value = 490.5
$ws.Range("M22").Value = -204
$ws.Range("N22").Value = -1080.5

$ws.Range("H27").Value = 493.33334
$ws.Range("I27").Value = 499
$ws.Range("J27").Value = 490.5
$ws.Range("K27").Value = 499
$ws.Range("L27").Value = 490.5
$ws.Range("M27").Value = -392
$ws.Range("N27").Value = -704.5

$ws.Range("H61").Value = 21830.4
$ws.Range("I61").Value = 26538
$ws.Range("K61").Value = 26538
$ws.Range("M61").Value = -26336

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H113").Value = 21830.4
$ws.Range("I113").Value = 26538
$ws.Range("K113").Value = 26538
$ws.Range("M113").Value = -24368

$ws.Range("H132").Value = 3121.1035
$ws.Range("I132").Value = 2475
$ws.Range("J132").Value = 4348.7
$ws.Range("K132").Value = 7425
$ws.Range("L132").Value = 13046.1
$ws.Range("M132").Value = -4895
$ws.Range("N132").Value = -18106.1

$ws.Range("H136").Value = 5530.143
$ws.Range("I136").Value = 3008.611
$ws.Range("J136").Value = 8200
$ws.Range("K136").Value = 9025.832999999999
$ws.Range("L136").Value = 24600
$ws.Range("M136").Value = -6475.832999999999
$ws.Range("N136").Value = -29700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 20690.375
$ws.Range("J74").Value = 24253.834
$ws.Range("L74").Value = 24253.834
$ws.Range("N74").Value = -26125.834

$ws.Range("H77").Value = 20690.375
$ws.Range("J77").Value = 24253.834
$ws.Range("L77").Value = 72761.50199999999
$ws.Range("N77").Value = -82121.50199999999

$ws.Range("H81").Value = 10003122
$ws.Range("I81").Value = 1826.3334
$ws.Range("J81").Value = 11768056
$ws.Range("K81").Value = 3652.6668
$ws.Range("L81").Value = 23536112
$ws.Range("M81").Value = -2591.6668
$ws.Range("N81").Value = -23538234

$ws.Range("H84").Value = 10003122
$ws.Range("I84").Value = 1826.3334
$ws.Range("J84").Value = 11768056
$ws.Range("K84").Value = 18263.334
$ws.Range("L84").Value = 117680560
$ws.Range("M84").Value = -12959.334
$ws.Range("N84").Value = -117691168

$ws.Range("H122").Value = 5634.2173
$ws.Range("I122").Value = 1352.4
$ws.Range("J122").Value = 13662.625
$ws.Range("K122").Value = 4057.2
$ws.Range("L122").Value = 40987.875
$ws.Range("M122").Value = -1607.2
$ws.Range("N122").Value = -45887.875

$ws.Range("H132").Value = 3062.0386
$ws.Range("I132").Value = 2832.3684
$ws.Range("J132").Value = 3685.4285
$ws.Range("K132").Value = 8497.1052
$ws.Range("L132").Value = 11056.2855
$ws.Range("M132").Value = -5967.1052
$ws.Range("N132").Value = -16116.2855
